$d = $word.ActiveDocument

$d.Content.Find.Execute("18÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "93÷8=", 2)
$d.Content.Find.Execute("62÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷9=", 2)
$d.Content.Find.Execute("85÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷3=", 2)
$d.Content.Find.Execute("80÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "73÷5=", 2)
$d.Content.Find.Execute("37÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "92÷2=", 2)
$d.Content.Find.Execute("98÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "99÷7=", 2)
$d.Content.Find.Execute("85÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "68÷5=", 2)
$d.Content.Find.Execute("75÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "38÷2=", 2)
$d.Content.Find.Execute("82÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "21÷5=", 2)
$d.Content.Find.Execute("69÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "75÷3=", 2)
$d.Content.Find.Execute("15÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "44÷6=", 2)
$d.Content.Find.Execute("52÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "41÷2=", 2)
$d.Content.Find.Execute("21÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "59÷4=", 2)
$d.Content.Find.Execute("74÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "28÷7=", 2)
$d.Content.Find.Execute("81÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "62÷3=", 2)
$d.Content.Find.Execute("39÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "24÷7=", 2)
$d.Content.Find.Execute("13÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "57÷6=", 2)
$d.Content.Find.Execute("23÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "39÷2=", 2)
$d.Content.Find.Execute("56÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷7=", 2)
$d.Content.Find.Execute("94÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "84÷4=", 2)
$d.Content.Find.Execute("17÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "64÷6=", 2)
$d.Content.Find.Execute("69÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "37÷8=", 2)
$d.Content.Find.Execute("93÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "22÷7=", 2)
$d.Content.Find.Execute("35÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "92÷4=", 2)
$d.Content.Find.Execute("28÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "29÷4=", 2)
